{"js": "// Append the \"chap6 and questions\" block to the end of the document body,\n// right before the final section break \u2014 a divider line, four Q&A\n// paragraphs, a closing divider line, and two trailing blank paragraphs.\n\nconst dashes1 = \"-\".repeat(135);\nconst dashes2 = \"-\".repeat(133);\n\n// OOXML fragment describing the new paragraphs, including the\n// lastRenderedPageBreak marker on the opening divider and the\n// proofErr spell-check bookmarks around the flagged words, to mirror\n// the authored markup exactly. A trailing empty paragraph is added to\n// the fragment to compensate for the paragraph mark that the\n// insertOoxml(..., replace) call below consumes from the placeholder\n// paragraph it is replacing.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:lastRenderedPageBreak/>\n              <w:t>${dashes1}</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">1. </w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\">I would describe it as the frame of an object. Using scaffold creates everything needed to for the frame to stand on its own, but there is nothing else added. In the case of rails it creates all the files needed for your application to run as well as files you are likely to use such as </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>css</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> sheets, but leaves them in a simple state that you will need to add to later. Scaffolding is traditionally used to mean a temporary structure used to support the process of building a structure, which is very similar to its role with rails.</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>2. In Rails, changes to a table definition are made through a migration file. Once we have made changes to this table in the migration file, we use the 'rake' command to apply the changes to the actual database. These changes can affect the schema of the database as well as the data in it.</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>3. test data is added</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\"> by editing the seeds file. Once the seeds file is populated you can populate the corresponding database with the data by running the rake command on the seed.</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">4. You find your </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>stylesheets</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> under app/assets/</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>stlyesheets</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t>. You do not need to run a separate command to generate them, the generate scaffold already created them.</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>${dashes2}</w:t>\n            </w:r>\n          </w:p>\n          <w:p/>\n          <w:p/>\n          <w:p/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n// Create a fresh empty paragraph at the very end of the body, then\n// replace it in place with the OOXML fragment above so the new\n// content lands after \"4. The controller manages sessions\" and before\n// the section break, exactly like the diff.\nconst body = context.document.body;\nconst placeholder = body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n\nplaceholder.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Append the \"chap6 and questions\" block to the end of the document \u2014\n# an opening divider line, four Q&A paragraphs, a closing divider line,\n# and two trailing blank paragraphs \u2014 right before the final section\n# break, mirroring the authored markup (including the proofErr\n# spell-check bookmarks and the lastRenderedPageBreak marker).\n\n$d = $word.ActiveDocument\n\n$dashes1 = \"---------------------------------------------------------------------------------------------------------------------------------------\"\n$dashes2 = \"-------------------------------------------------------------------------------------------------------------------------------------\"\n\n$wns = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\n$xml = @\"\n<w:p $wns><w:r><w:lastRenderedPageBreak/><w:t>$dashes1</w:t></w:r></w:p><w:p $wns><w:r><w:t xml:space=\"preserve\">1. </w:t></w:r><w:r $wns><w:t xml:space=\"preserve\">I would describe it as the frame of an object. Using scaffold creates everything needed to for the frame to stand on its own, but there is nothing else added. In the case of rails it creates all the files needed for your application to run as well as files you are likely to use such as </w:t></w:r><w:proofErr $wns w:type=\"spellStart\"/><w:r $wns><w:t>css</w:t></w:r><w:proofErr $wns w:type=\"spellEnd\"/><w:r $wns><w:t xml:space=\"preserve\"> sheets, but leaves them in a simple state that you will need to add to later. Scaffolding is traditionally used to mean a temporary structure used to support the process of building a structure, which is very similar to its role with rails.</w:t></w:r></w:p><w:p $wns><w:r><w:t>2. In Rails, changes to a table definition are made through a migration file. Once we have made changes to this table in the migration file, we use the 'rake' command to apply the changes to the actual database. These changes can affect the schema of the database as well as the data in it.</w:t></w:r></w:p><w:p $wns><w:r><w:t>3. test data is added</w:t></w:r><w:r $wns><w:t xml:space=\"preserve\"> by editing the seeds file. Once the seeds file is populated you can populate the corresponding database with the data by running the rake command on the seed.</w:t></w:r></w:p><w:p $wns><w:r><w:t xml:space=\"preserve\">4. You find your </w:t></w:r><w:proofErr $wns w:type=\"spellStart\"/><w:r $wns><w:t>stylesheets</w:t></w:r><w:proofErr $wns w:type=\"spellEnd\"/><w:r $wns><w:t xml:space=\"preserve\"> under app/assets/</w:t></w:r><w:proofErr $wns w:type=\"spellStart\"/><w:r $wns><w:t>stlyesheets</w:t></w:r><w:proofErr $wns w:type=\"spellEnd\"/><w:r $wns><w:t>. You do not need to run a separate command to generate them, the generate scaffold already created them.</w:t></w:r></w:p><w:p $wns><w:r><w:t>$dashes2</w:t></w:r></w:p><w:p $wns/><w:p $wns/>\n\"@\n\n$endRange = $d.Content\n$endRange.Collapse(0)\n[void]$endRange.InsertXML($xml)\n"}
